$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 213 (Serie 227 - 25-10-2021) with revised figures ---
$ws.Range("B213").Value = 1.63
$ws.Range("C213").Value = 1.14
$ws.Range("E213").Value = -0.11
$ws.Range("J213").Value = 8.03
$ws.Range("L213").Value = 7.12
$ws.Range("M213").Value = 12.03
$ws.Range("N213").Value = 7.9
$ws.Range("O213").Value = 7.56

# --- Row 214 (Serie 228 - 26-10-2021) ---
$ws.Range("A214").Value = "26-10-2021"
$ws.Range("B214").Value = 1.61
$ws.Range("C214").Value = 1.11
$ws.Range("D214").Value = 0.11
$ws.Range("E214").Value = -0.12
$ws.Range("F214").Value = 2.49
$ws.Range("G214").Value = 3.59
$ws.Range("H214").Value = 2.87
$ws.Range("I214").Value = 2.56
$ws.Range("J214").Value = 8.08
$ws.Range("K214").Value = 1.96
$ws.Range("L214").Value = 6.96
$ws.Range("M214").Value = 12
$ws.Range("O214").Value = 7.57
$ws.Range("P214").Value = 5.94

# --- Row 215 (Serie 229 - 27-10-2021) ---
$ws.Range("A215").Value = "27-10-2021"
$ws.Range("B215").Value = 1.54
$ws.Range("C215").Value = 0.98
$ws.Range("D215").Value = 0.1
$ws.Range("E215").Value = -0.18
$ws.Range("F215").Value = 2.55
$ws.Range("G215").Value = 3.58
$ws.Range("H215").Value = 2.81
$ws.Range("I215").Value = 2.58
$ws.Range("J215").Value = 8.04
$ws.Range("K215").Value = 1.94
$ws.Range("L215").Value = 6.87
$ws.Range("M215").Value = 11.9
$ws.Range("N215").Value = 8.04
$ws.Range("O215").Value = 7.53
$ws.Range("P215").Value = 5.97

# --- Row 216 (Serie 230 - 28-10-2021) ---
$ws.Range("A216").Value = "28-10-2021"
$ws.Range("B216").Value = 1.58
$ws.Range("C216").Value = 1.01
$ws.Range("D216").Value = 0.09
$ws.Range("E216").Value = -0.14
$ws.Range("F216").Value = 2.46
$ws.Range("G216").Value = 3.57
$ws.Range("H216").Value = 2.85
$ws.Range("J216").Value = 8.19
$ws.Range("K216").Value = 1.92
$ws.Range("L216").Value = 6.89
$ws.Range("M216").Value = 12.42
$ws.Range("O216").Value = 7.54
$ws.Range("P216").Value = 5.82

# --- Row 217 (Serie 231 - 29-10-2021) ---
$ws.Range("A217").Value = "29-10-2021"
$ws.Range("B217").Value = 1.56
$ws.Range("C217").Value = 1.03
$ws.Range("D217").Value = 0.1
$ws.Range("E217").Value = -0.11
$ws.Range("F217").Value = 2.6
$ws.Range("G217").Value = 3.6
$ws.Range("H217").Value = 2.82
$ws.Range("I217").Value = 2.6
$ws.Range("J217").Value = 8.22
$ws.Range("K217").Value = 1.98
$ws.Range("L217").Value = 6.92
$ws.Range("M217").Value = 12.22
$ws.Range("O217").Value = 7.49
$ws.Range("P217").Value = 5.83

# --- Row 218 (Serie 232 - 01-11-2021) ---
# "01-11-2021" looks like an ambiguous M-D-Y date to the engine's autodetection,
# so force the cell to Text first, then drop the explicit style again so the
# saved cell has no number format / style applied (matches a plain shared string).
$ws.Range("A218").NumberFormat = "@"
$ws.Range("A218").Value = "01-11-2021"
$ws.Range("A218").Style = "Normal"
$ws.Range("B218").Value = 1.56
$ws.Range("C218").Value = 1.06
$ws.Range("D218").Value = 0.1
$ws.Range("E218").Value = -0.1
$ws.Range("F218").Value = 2.54
$ws.Range("G218").Value = 3.63
$ws.Range("I218").Value = 2.59
$ws.Range("J218").Value = 8.35
$ws.Range("K218").Value = 1.99
$ws.Range("L218").Value = 6.86
$ws.Range("M218").Value = 12.35
$ws.Range("O218").Value = 7.59

# --- Row 219 (Serie 233 - 02-11-2021) ---
$ws.Range("A219").NumberFormat = "@"
$ws.Range("A219").Value = "02-11-2021"
$ws.Range("A219").Style = "Normal"
$ws.Range("B219").Value = 1.56
$ws.Range("C219").Value = 1.05
$ws.Range("D219").Value = 0.08
$ws.Range("E219").Value = -0.15
$ws.Range("F219").Value = 2.49
$ws.Range("G219").Value = 3.63
$ws.Range("H219").Value = 2.85
$ws.Range("I219").Value = 2.56
$ws.Range("J219").Value = 8.31
$ws.Range("K219").Value = 1.95
$ws.Range("L219").Value = 6.72
